$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4 (G4=5470)
$ws.Range("H4").Value = 236.25
$ws.Range("I4").Value = 245
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 245
$ws.Range("L4").Value = 140
$ws.Range("M4").Value = -131
$ws.Range("N4").Value = -368

# Row 6 (G6=4564)
$ws.Range("H6").Value = 156.35715
$ws.Range("I6").Value = 143.76923
$ws.Range("K6").Value = 431.30769
$ws.Range("M6").Value = -319.30769

# Row 70 (G70=12604)
$ws.Range("H70").Value = 1652.8182
$ws.Range("I70").Value = 1532.8572
$ws.Range("J70").Value = 1862.75
$ws.Range("K70").Value = 4598.571599999999
$ws.Range("L70").Value = 5588.25
$ws.Range("M70").Value = -4328.571599999999
$ws.Range("N70").Value = -6128.25

# Row 73 (G73=12604)
$ws.Range("H73").Value = 1652.8182
$ws.Range("I73").Value = 1532.8572
$ws.Range("J73").Value = 1862.75
$ws.Range("K73").Value = 4598.571599999999
$ws.Range("L73").Value = 5588.25
$ws.Range("M73").Value = -3662.571599999999
$ws.Range("N73").Value = -7460.25

# Row 113 (G113=27775)
$ws.Range("H113").Value = 8716.821
$ws.Range("J113").Value = 8193.471
$ws.Range("L113").Value = 8193.471
$ws.Range("N113").Value = -14701.471

# Row 137 (G137=44013)
$ws.Range("H137").Value = 2969.1875
$ws.Range("I137").Value = 1552.0714
$ws.Range("K137").Value = 4656.2142
$ws.Range("M137").Value = -2106.2142

$ws = $wb.Worksheets.Item("ARM")
# Row 4 (G4=5071)
$ws.Range("H4").Value = 229.66667
$ws.Range("I4").Value = 245
$ws.Range("K4").Value = 245
$ws.Range("M4").Value = -129

# Row 5 (G5=5091)
$ws.Range("H5").Value = 369.75
$ws.Range("I5").Value = 369.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 369.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -257.75
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 4 (G4=5091)
$ws.Range("H4").Value = 369.75
$ws.Range("I4").Value = 369.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 369.75
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -254.75
$ws.Range("N4").ClearContents()

# Row 22 (G22=5092)
$ws.Range("H22").Value = 542.9231
$ws.Range("I22").Value = 564.36365
$ws.Range("J22").Value = 425
$ws.Range("K22").Value = 564.36365
$ws.Range("L22").Value = 425
$ws.Range("M22").Value = -391.36365
$ws.Range("N22").Value = -771

# Row 25 (G25=2370)
$ws.Range("H25").Value = 883.2
$ws.Range("I25").Value = 800
$ws.Range("K25").Value = 800
$ws.Range("M25").Value = -565

# Row 86 (G86=12526)
$ws.Range("H86").Value = 2623.7827
$ws.Range("I86").Value = 2560.7334
$ws.Range("K86").Value = 2560.7334
$ws.Range("M86").Value = -1437.7334

# Row 89 (G89=12526)
$ws.Range("H89").Value = 2623.7827
$ws.Range("I89").Value = 2560.7334
$ws.Range("K89").Value = 12803.667
$ws.Range("M89").Value = -7187.667000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (G4=3742)
$ws.Range("H4").Value = 4522.6924
$ws.Range("I4").Value = 3691.25
$ws.Range("J4").Value = 14500
$ws.Range("K4").Value = 3691.25
$ws.Range("L4").Value = 14500
$ws.Range("M4").Value = -3579.25
$ws.Range("N4").Value = -14724

# Row 7 (G7=5361)
$ws.Range("H7").Value = 98.611115
$ws.Range("J7").Value = 121.666664
$ws.Range("L7").Value = 121.666664
$ws.Range("N7").Value = -347.666664

# Row 16 (G16=27691)
$ws.Range("H16").Value = 1018.61536
$ws.Range("I16").Value = 723.8333
$ws.Range("K16").Value = 723.8333
$ws.Range("M16").Value = -436.8333

# Row 31 (G31=44023)
$ws.Range("H31").Value = 2230.6858
$ws.Range("I31").Value = 2156.28
$ws.Range("J31").Value = 2416.7
$ws.Range("K31").Value = 2156.28
$ws.Range("L31").Value = 2416.7
$ws.Range("M31").Value = -1861.28
$ws.Range("N31").Value = -3006.7

# Row 34 (G34=44023)
$ws.Range("H34").Value = 2230.6858
$ws.Range("I34").Value = 2156.28
$ws.Range("J34").Value = 2416.7
$ws.Range("K34").Value = 2156.28
$ws.Range("L34").Value = 2416.7
$ws.Range("M34").Value = -1954.28
$ws.Range("N34").Value = -2820.7

# Row 113 (G113=27691)
$ws.Range("H113").Value = 1018.61536
$ws.Range("I113").Value = 723.8333
$ws.Range("K113").Value = 723.8333
$ws.Range("M113").Value = 1446.1667

# Row 122 (G122=36196)
$ws.Range("H122").Value = 2473.6072
$ws.Range("I122").Value = 1555.0625
$ws.Range("J122").Value = 3698.3333
$ws.Range("K122").Value = 4665.1875
$ws.Range("L122").Value = 11094.9999
$ws.Range("M122").Value = -2215.1875
$ws.Range("N122").Value = -15994.9999

# Row 132 (G132=44019)
$ws.Range("H132").Value = 1362
$ws.Range("I132").Value = 1270.125
$ws.Range("K132").Value = 3810.375
$ws.Range("M132").Value = -1280.375

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (G2=4847)
$ws.Range("H2").Value = 16666960
$ws.Range("J2").Value = 41667230
$ws.Range("L2").Value = 250003380
$ws.Range("N2").Value = -250003606

# Row 4 (G4=4650)
$ws.Range("H4").Value = 1642.25
$ws.Range("I4").Value = 2139.6667
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 6419.000100000001
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = -6307.000100000001
$ws.Range("N4").Value = -674

# Row 6 (G6=4639)
$ws.Range("H6").Value = 1180
$ws.Range("I6").Value = 225
$ws.Range("K6").Value = 675
$ws.Range("M6").Value = -562

# Row 10 (G10=4689)
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

# Row 16 (G16=4641)
$ws.Range("H16").Value = 168
$ws.Range("I16").Value = 114.666664
$ws.Range("K16").Value = 343.999992
$ws.Range("M16").Value = -170.999992

# Row 131 (G131=36060)
$ws.Range("H131").Value = 2338.6667
$ws.Range("I131").Value = 1012.5
$ws.Range("J131").Value = 3399.6
$ws.Range("K131").Value = 3037.5
$ws.Range("L131").Value = 10198.8
$ws.Range("M131").Value = 2002.5
$ws.Range("N131").Value = -20278.8

# Row 132 (G132=43972)
$ws.Range("H132").Value = 1381.625
$ws.Range("I132").Value = 1263.25
$ws.Range("K132").Value = 11369.25
$ws.Range("M132").Value = -8839.25

# Row 137 (G137=44088)
$ws.Range("H137").Value = 2642.2222
$ws.Range("I137").Value = 2533.1667
$ws.Range("J137").Value = 2860.3333
$ws.Range("K137").Value = 7599.500100000001
$ws.Range("L137").Value = 8580.999899999999
$ws.Range("M137").Value = -2499.500100000001
$ws.Range("N137").Value = -18780.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (G2=5062)
$ws.Range("H2").Value = 283.1111
$ws.Range("I2").Value = 270.83334
$ws.Range("K2").Value = 270.83334
$ws.Range("M2").Value = -157.83334

$ws = $wb.Worksheets.Item("LTW")
# Row 68 (G68=12563)
$ws.Range("H68").Value = 5903.8184
$ws.Range("J68").Value = 11498.5
$ws.Range("L68").Value = 11498.5
$ws.Range("N68").Value = -12996.5

# Row 71 (G71=12563)
$ws.Range("H71").Value = 5903.8184
$ws.Range("J71").Value = 11498.5
$ws.Range("L71").Value = 57492.5
$ws.Range("N71").Value = -64980.5

# Row 122 (G122=36247)
$ws.Range("H122").Value = 8934.615
$ws.Range("I122").Value = 8540.049999999999
$ws.Range("K122").Value = 25620.15
$ws.Range("M122").Value = -23170.15

# Row 132 (G132=44058)
$ws.Range("H132").Value = 4294.569
$ws.Range("I132").Value = 4158.891
$ws.Range("K132").Value = 12476.673
$ws.Range("M132").Value = -9946.672999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (G2=3307)
$ws.Range("H2").Value = 29001.8
$ws.Range("I2").Value = 5002
$ws.Range("J2").Value = 35001.75
$ws.Range("K2").Value = 5002
$ws.Range("L2").Value = 35001.75
$ws.Range("M2").Value = -4890
$ws.Range("N2").Value = -35225.75

# Row 68 (G68=10762)
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71 (G71=10762)
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 132 (G132=44029)
$ws.Range("H132").Value = 2075.0715
$ws.Range("I132").Value = 2050.077
$ws.Range("K132").Value = 6150.231000000001
$ws.Range("M132").Value = -3620.231000000001
